$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting the existing row 9 (and below) down to row 10.
$ws.Rows.Item(9).Insert()

# The inserted row 9 copied formatting from the row above (row 8); now fill it with the
# data that used to live in the old row 9 (now row 10 after the shift).
$ws.Range("A9").Value = 12
$ws.Range("B9").Value = "Mapocho Venta Directa de Santiago"
$ws.Range("C9").Value = "Metropolitana"
$ws.Range("D9").Value = 44435
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E9").Value = 13
$ws.Range("F9").Value = 100112013
$ws.Range("G9").Value = "Alcachofa"
$ws.Range("H9").Value = "Española"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 25
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 14000
$ws.Range("N9").Value = "`$/caja 30 unidades"
$ws.Range("O9").Value = "Provincia del Elquí"
$ws.Range("P9").Value = 467
$ws.Range("Q9").Value = 30
$ws.Range("R9").Value = "Hortaliza"

# Row 7 now contains new data values.
$ws.Range("D7").Value = 44449
$ws.Range("J7").Value = 45
$ws.Range("K7").Value = 12000
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = 12000
$ws.Range("P7").Value = 400

# Row 8: Origen changes.
$ws.Range("O8").Value = "Provincia de Limarí"
